$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 115
$ws1.Range("F5").Value = 3022
$ws1.Range("F6").Value = 305
$ws1.Range("F7").Value = 406

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 115
$ws4.Range("F5").Value = 3022
$ws4.Range("F6").Value = 305
$ws4.Range("F9").Value = 406
